# feat: add 2022-Q1 data
#
# The original "总计" (summary) sheet becomes the new "2022-Q1" sheet
# (repurposed in place, keeping its original sheetId) populated with the
# 2022-Q1 fund-holding detail table, and a brand-new "总计" sheet is
# appended after it containing the updated summary table (with the new
# 2022-Q1 row prepended).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the current "总计" sheet BEFORE we touch it, so the
# duplicate inherits all of its sheet-level formatting (sheetPr,
# pageMargins, sheetFormatPr, column-A number style, etc.). The copy is
# placed immediately after the source sheet, which is exactly the final
# tab order we want ( ... , 2022-Q1, 总计 ).
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Copy([System.Type]::Missing, $oldTotal)
$wsTotal = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTotal.Name = "总计新"

# Give the brand-new last row (row 7) the same style as the other index
# cells in column A (s=2) before we populate it, by copying a cell that
# already carries that style.
$wsTotal.Range("A2").Copy($wsTotal.Range("A7"))

# ---------------------------------------------------------------------
# Step 2: rewrite the duplicate with the updated summary table: the new
# 2022-Q1 row on top, followed by the previously-existing rows (each
# index bumped by one).
# ---------------------------------------------------------------------
$summaryRows = @(
    @("2022-Q1", 5, 2.56),
    @("2021-Q4", 10, 3.25),
    @("2021-Q3", 12, 4.49),
    @("2021-Q2", 22, 3.78),
    @("2021-Q1", 17, 2.73),
    @("2020-Q4", 23, 7.1)
)

$wsTotal.Cells.Item(1, 2).Value = "日期"
$wsTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$wsTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

$r = 2
foreach ($row in $summaryRows) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
    $wsTotal.Cells.Item($r, 2).Value = $row[0]
    $wsTotal.Cells.Item($r, 3).Value = $row[1]
    $wsTotal.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 3: turn the original "总计" worksheet object into "2022-Q1": seed
# its formatting from the "2021-Q4" detail sheet (same 8-column detail
# layout), then overwrite with the 2022-Q1 fund-holding data.
# ---------------------------------------------------------------------
$ws2022 = $oldTotal
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H11").Copy($ws2022.Range("A1"))
$ws2022.Range("A1").Clear()
$ws2022.Range("A7:H11").Clear()
$ws2022.Name = "2022-Q1"

# Force text storage for the columns whose values must not be
# reinterpreted as numbers (leading zeros in fund codes, fixed decimal
# display for the percentage/size/value columns).
$ws2022.Range("B2:B6").NumberFormat = "@"
$ws2022.Range("D2:G6").NumberFormat = "@"

$ws2022.Cells.Item(1, 2).Value = "基金代码"
$ws2022.Cells.Item(1, 3).Value = "基金名称"
$ws2022.Cells.Item(1, 4).Value = "基金规模"
$ws2022.Cells.Item(1, 5).Value = "股票总仓位"
$ws2022.Cells.Item(1, 6).Value = "仓位占比"
$ws2022.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws2022.Cells.Item(1, 8).Value = "仓位排名"

$fundRows = @(
    @("161914", "万家创业板2年定期开放混合A", "15.74", "95.20", "6.49", "1.0215", 8),
    @("070001", "嘉实成长收益混合A", "24.07", "72.54", "3.30", "0.7943", 8),
    @("160916", "大成优选混合(LOF)", "16.14", "89.35", "3.67", "0.5923", 9),
    @("161915", "万家创业板2年定期开放混合C", "2.36", "95.20", "6.49", "0.1532", 8),
    @("960024", "嘉实成长收益混合H", "0.01", "72.54", "3.30", "0.0003", 8)
)

$r = 2
foreach ($row in $fundRows) {
    $ws2022.Cells.Item($r, 1).Value = $r - 2
    $ws2022.Cells.Item($r, 2).Value = $row[0]
    $ws2022.Cells.Item($r, 3).Value = $row[1]
    $ws2022.Cells.Item($r, 4).Value = $row[2]
    $ws2022.Cells.Item($r, 5).Value = $row[3]
    $ws2022.Cells.Item($r, 6).Value = $row[4]
    $ws2022.Cells.Item($r, 7).Value = $row[5]
    $ws2022.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 4: now that the original "总计" name is freed up (it became
# "2022-Q1"), give the new summary sheet its proper final name.
# ---------------------------------------------------------------------
$wsTotal.Name = "总计"
